# Add a "popularity" column (column S) to the top2017 sheet, matching the
# Spotify popularity score scraped for each track.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell in S1, formatted like the other header cells (bold, thin box
# border, centered horizontally, top-aligned vertically).
$ws.Range("S1").Value = "popularity"
$ws.Range("S1").Font.Bold = $true
$ws.Range("S1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("S1").VerticalAlignment = -4160     # xlTop
$ws.Range("S1").Borders.LineStyle = 1         # xlContinuous (thin box border)

# Popularity values for rows 2..101 (same song order as the existing data).
$popularity = @(
    87,76,80,84,83,85,69,82,83,70,
    85,64,14,83,82,84,77,79,78,89,
    74,86,86,83,78,72,81,82,86,78,
    81,85,68,80,84,82,81,89,76,0,
    89,68,72,74,78,78,83,76,77,85,
    59,76,77,81,83,75,78,81,82,79,
    80,72,53,72,67,67,39,78,77,77,
    81,87,84,75,66,82,77,62,75,77,
    77,75,71,64,75,84,81,70,78,88,
    62,77,74,76,81,75,81,78,72,82
)

for ($i = 0; $i -lt $popularity.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 19).Value = $popularity[$i]
}
